# Applies the "#5: insurance, claim, debt, investment done" edit:
# fixes the header row (row 1) of the 保險 (insurance) and 債務 (debt) sheets,
# which incorrectly duplicated the first data row, and appends the standard
# trailing metadata columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) that every other
# sheet in this workbook already has.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "保險" (insurance) -- currently only columns B:D, rows 1:2
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item("保險")

# Copy the existing header/data formatting onto the new trailing columns
# before filling them in, so E:K pick up the same cell styles as B:D.
$wsIns.Range("B1:D1").Copy() | Out-Null
$wsIns.Range("E1:K1").PasteSpecial(-4122) | Out-Null
$wsIns.Range("B2:D2").Copy() | Out-Null
$wsIns.Range("E2:K2").PasteSpecial(-4122) | Out-Null

# Row 1: proper column headers (previously this row wrongly repeated the
# row-2 data values instead of naming the fields).
$wsIns.Range("B1").Value = "company"
$wsIns.Range("C1").Value = "name"
$wsIns.Range("D1").Value = "owner"
$wsIns.Range("E1").Value = "property_category"
$wsIns.Range("F1").Value = "category"
$wsIns.Range("G1").Value = "date"
$wsIns.Range("H1").Value = "legislator_name"
$wsIns.Range("I1").Value = "legislator_id"
$wsIns.Range("J1").Value = "source_file"
$wsIns.Range("K1").Value = "index"

# Row 2: B:D already hold the correct data (company / name / owner) --
# append the standard metadata columns.
$wsIns.Range("E2").Value = "insurance"
$wsIns.Range("F2").Value = "normal"
$wsIns.Range("G2").Value = "2011-12-06"
$wsIns.Range("H2").Value = "王廷升"
$wsIns.Range("I2").Value = 1727
$wsIns.Range("J2").Value = "tmp44311"
$wsIns.Range("K2").Value = 97

# ---------------------------------------------------------------------
# Sheet "債務" (debt) -- currently only columns B:G, rows 1:4
# ---------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item("債務")

$wsDebt.Range("B1:G1").Copy() | Out-Null
$wsDebt.Range("H1:N1").PasteSpecial(-4122) | Out-Null
$wsDebt.Range("B2:G2").Copy() | Out-Null
$wsDebt.Range("H2:N4").PasteSpecial(-4122) | Out-Null

# Row 1: proper column headers (previously wrongly repeated row-2 data).
$wsDebt.Range("B1").Value = "species"
$wsDebt.Range("C1").Value = "debtor"
$wsDebt.Range("D1").Value = "owner"
$wsDebt.Range("E1").Value = "total"
$wsDebt.Range("F1").Value = "register_date"
$wsDebt.Range("G1").Value = "register_reason"
$wsDebt.Range("H1").Value = "property_category"
$wsDebt.Range("I1").Value = "category"
$wsDebt.Range("J1").Value = "date"
$wsDebt.Range("K1").Value = "legislator_name"
$wsDebt.Range("L1").Value = "legislator_id"
$wsDebt.Range("M1").Value = "source_file"
$wsDebt.Range("N1").Value = "index"

# Rows 2-4: B:G already hold the correct data -- append metadata columns.
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").Value = "2011-12-06"
$wsDebt.Range("K2").Value = "王廷升"
$wsDebt.Range("L2").Value = 1727
$wsDebt.Range("M2").Value = "tmp44311"
$wsDebt.Range("N2").Value = 112

$wsDebt.Range("H3").Value = "debt"
$wsDebt.Range("I3").Value = "normal"
$wsDebt.Range("J3").Value = "2011-12-06"
$wsDebt.Range("K3").Value = "王廷升"
$wsDebt.Range("L3").Value = 1727
$wsDebt.Range("M3").Value = "tmp44311"
$wsDebt.Range("N3").Value = 113

$wsDebt.Range("H4").Value = "debt"
$wsDebt.Range("I4").Value = "normal"
$wsDebt.Range("J4").Value = "2011-12-06"
$wsDebt.Range("K4").Value = "王廷升"
$wsDebt.Range("L4").Value = 1727
$wsDebt.Range("M4").Value = "tmp44311"
$wsDebt.Range("N4").Value = 114
